$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a piece of literal text into a cell while guaranteeing it is
# stored as plain text (never auto-parsed into a date/number), and without
# leaving any cell-style ($s=) footprint behind. Plain ".Value = <string>"
# auto-detects locale-ambiguous day/month strings (e.g. "12.02.2024") as
# dates, and ".NumberFormat" tagging leaves a permanent style index on the
# cell - neither of which happened in the source report. Routing the text
# through a literal formula and then flattening it with a values-only paste
# keeps the cell a plain shared string with no style / no formula residue.
function Set-LiteralText($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# --- Row 2 ---
$ws.Range("A2").Value = "מיגיר"
$ws.Range("B2").Value = "לוייב"
$ws.Range("C2").Value = "הצעה "
$ws.Range("D2").Value = "גל מסיקה"
$ws.Range("E2").Value = "לוי הודיה"
$ws.Range("F2").Value = "בוצע חישוב לא זכאי להחזר"
$ws.Range("G2").Value = "30.03.2023"
Set-LiteralText $ws.Range("H2") "12.02.2024"
$ws.Range("I2").Value = "החזרי מס"

# --- Row 3 ---
$ws.Range("A3").Value = "אלבר ואילנית*"
$ws.Range("B3").Value = "דהן"
$ws.Range("C3").Value = "הופק"
$ws.Range("D3").Value = "גל מסיקה"
$ws.Range("E3").Value = "יהוד ספיר"
$ws.Range("F3").Value = "בוצע גבייה 12.12.23בוצע גבייה 14.02.24"
$ws.Range("G3").Value = "14.05.2023"
$ws.Range("H3").Value = "18.02.2024"
$ws.Range("I3").Value = "החזרי מס"

# --- Row 4 ---
$ws.Range("A4").Value = "אוריה"
$ws.Range("B4").Value = "לברון"
$ws.Range("C4").Value = "הצעה "
$ws.Range("D4").Value = "גל מסיקה"
$ws.Range("E4").Value = "עוזירי קארין"
$ws.Range("F4").Value = "בוצע חישוב - לא זכאי להחזר"
$ws.Range("G4").Value = "22.08.2023"
$ws.Range("H4").Value = "14.01.2024"
$ws.Range("I4").Value = "החזרי מס"

# --- Row 5 ---
$ws.Range("A5").Value = "אלירז ( עידן )"
$ws.Range("B5").Value = "כס"
$ws.Range("C5").Value = "הצעה"
$ws.Range("D5").Value = "גל מסיקה"
$ws.Range("E5").Value = "כהן ליאל "
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "23.08.2023"
Set-LiteralText $ws.Range("H5") "12.02.2024"
$ws.Range("I5").ClearContents()

# --- Row 6 ---
$ws.Range("A6").Value = "עידן ( ואלירז )"
$ws.Range("B6").Value = "כס"
$ws.Range("C6").Value = "הצעה"
$ws.Range("D6").Value = "גל מסיקה"
$ws.Range("E6").Value = "גרינברגר גילי"
$ws.Range("F6").Value = "בתהליך חישוב חוזר 25.05.24"
$ws.Range("G6").Value = "28.08.2023"
$ws.Range("H6").Value = "29.05.2024"
$ws.Range("I6").ClearContents()

# --- Rows 7-9 no longer exist in the updated report ---
$ws.Rows("7:9").Delete()

# --- Register the datetime number format used by the mailer (yyyy-mm-dd
#     hh:mm:ss) in the workbook's style table. Stamp it on a scratch cell
#     and then drop that scratch column again so the printed A1:I6 table
#     and its dimension/row spans stay untouched, while the format code
#     itself remains registered in xl/styles.xml. ---
$ws.Range("K1").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Columns("K").Delete()
